# Update the cryptos price/volume table to the latest scrape.
# Each row corresponds to one coin; columns are:
#   B = Coin name, C = Link, D = Price, E = Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.912.19"

$ws.Cells.Item(3, 4).Value = "3.375.71"
$ws.Cells.Item(3, 5).Value = "  -4.49%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).Value = "561.45"
$ws.Cells.Item(5, 5).Value = "  -4.13%  "

$ws.Cells.Item(6, 4).Value = "184.99"
$ws.Cells.Item(6, 5).Value = "  -5.95%  "

$ws.Cells.Item(7, 4).Value = "0.599"
$ws.Cells.Item(7, 5).Value = "  -2.17%  "

$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(9, 4).Value = "3.362.70"
$ws.Cells.Item(9, 5).Value = "  -4.50%  "

$ws.Cells.Item(10, 5).Value = "  -8.46%  "

$ws.Cells.Item(11, 4).Value = "0.599"
$ws.Cells.Item(11, 5).Value = "  -4.70%  "

$ws.Cells.Item(12, 4).Value = "48.09"
$ws.Cells.Item(12, 5).Value = "  -7.36%  "

$ws.Cells.Item(13, 4).Value = "0.0000272"
$ws.Cells.Item(13, 5).Value = "  -5.89%  "

$ws.Cells.Item(14, 4).Value = "8.79"
$ws.Cells.Item(14, 5).Value = "  -5.84%  "

$ws.Cells.Item(15, 4).Value = "3.909.28"
$ws.Cells.Item(15, 5).Value = "  -4.49%  "

$ws.Cells.Item(16, 4).Value = "605.28"
$ws.Cells.Item(16, 5).Value = "  -11.56%  "

$ws.Cells.Item(17, 4).Value = "66.724.25"
$ws.Cells.Item(17, 5).Value = "  -3.90%  "

$ws.Cells.Item(18, 4).Value = "3.370.51"
$ws.Cells.Item(18, 5).Value = "  -4.77%  "

$ws.Cells.Item(19, 4).Value = "18.02"
$ws.Cells.Item(19, 5).Value = "  -2.99%  "

$ws.Cells.Item(20, 5).Value = "  -2.85%  "

$ws.Cells.Item(21, 4).Value = "'11.70"
$ws.Cells.Item(21, 5).Value = "  -6.00%  "

$ws.Cells.Item(22, 4).Value = "0.918"
$ws.Cells.Item(22, 5).Value = "  -5.25%  "

$ws.Cells.Item(23, 5).Value = "  -4.52%  "

$ws.Cells.Item(24, 4).Value = "5.16"
$ws.Cells.Item(24, 5).Value = "  -1.60%  "

$ws.Cells.Item(25, 4).Value = "98.78"
$ws.Cells.Item(25, 5).Value = "  -8.64%  "

$ws.Cells.Item(26, 5).Value = "  -6.88%  "

$ws.Cells.Item(27, 5).Value = "  -0.12%  "

$ws.Cells.Item(28, 5).Value = "  -5.85%  "

$ws.Cells.Item(29, 4).Value = "9.56"
$ws.Cells.Item(29, 5).Value = "  -7.66%  "

$ws.Cells.Item(30, 5).Value = "  -8.87%  "

$ws.Cells.Item(31, 4).Value = "30.92"
$ws.Cells.Item(31, 5).Value = "  -7.56%  "

$ws.Cells.Item(32, 4).Value = "3.99"
$ws.Cells.Item(32, 5).Value = "  -9.46%  "

$ws.Cells.Item(33, 4).Value = "6.37"
$ws.Cells.Item(33, 5).Value = "  -7.97%  "

$ws.Cells.Item(34, 4).Value = "11.23"
$ws.Cells.Item(34, 5).Value = "  -5.84%  "

$ws.Cells.Item(35, 2).Value = "Bittensor"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(35, 4).Value = "549.88"
$ws.Cells.Item(35, 5).Value = "  +8.59%  "

$ws.Cells.Item(36, 2).Value = "Maker"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(36, 4).Value = "3.874.52"
$ws.Cells.Item(36, 5).Value = "  +1.92%  "

$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).Value = "0.106"
$ws.Cells.Item(37, 5).Value = "  -4.92%  "

$ws.Cells.Item(38, 4).Value = "58.47"
$ws.Cells.Item(38, 5).Value = "  -6.12%  "

$ws.Cells.Item(39, 5).Value = "  +0.01%  "

$ws.Cells.Item(40, 4).Value = "3.44"
$ws.Cells.Item(40, 5).Value = "  -4.32%  "

$ws.Cells.Item(41, 4).Value = "0.0₃0728"
$ws.Cells.Item(41, 5).Value = "  -11.20%  "

$ws.Cells.Item(42, 2).Value = "Fetch.AI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(42, 4).Value = "2.74"
$ws.Cells.Item(42, 5).Value = "  -7.75%  "

$ws.Cells.Item(43, 2).Value = "Kaspa"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43, 4).Value = "0.129"
$ws.Cells.Item(43, 5).Value = "  -5.31%  "

$ws.Cells.Item(44, 2).Value = "CoreDAO"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Cells.Item(44, 4).Value = "3.39"
$ws.Cells.Item(44, 5).Value = "  +26.12%  "

$ws.Cells.Item(45, 5).Value = "  -5.52%  "

$ws.Cells.Item(46, 4).Value = "32.63"
$ws.Cells.Item(46, 5).Value = "  -6.66%  "

$ws.Cells.Item(47, 4).Value = "0.0421"
$ws.Cells.Item(47, 5).Value = "  -7.89%  "

$ws.Cells.Item(48, 4).Value = "3.23"
$ws.Cells.Item(48, 5).Value = "  -4.33%  "

$ws.Cells.Item(49, 4).Value = "2.69"
$ws.Cells.Item(49, 5).Value = "  -9.17%  "

$ws.Cells.Item(50, 5).Value = "  -4.66%  "

$ws.Cells.Item(51, 5).Value = "  -0.23%  "
